# Update weekly Fruta/Hortaliza price data for "Vega Modelo de Temuco - Chirimoya"
# Rows 12-28: update existing rows (date, calidad, volumen, precios, unidad, precio/kg, kg/unidad)
# Rows 29-36: replace/append rows with full data (new week's records)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: rows 12-28 already have the constant columns (A,B,C,E,F,G,H,I,J,K,R) correct;
#     only D (Fecha), L (Calidad), M (Volumen), N (Precio minimo), O (Precio maximo),
#     P (Precio promedio ponderado), Q (Unidad de comercializacion), S (Precio $/Kg) and
#     T (Kg / unidad) change.
$existingRows = @(
    @(12, 44434, "Especial", 80,  3500,  3500,  3500,  "$/kilo (en caja de 15 kilos)", 3500, 1),
    @(13, 44428, "Primera",  55,  3500,  3500,  3500,  "$/kilo (en caja de 15 kilos)", 3500, 1),
    @(14, 44162, "Primera",  85,  2200,  2300,  2247,  "$/kilo (en caja de 15 kilos)", 2247, 1),
    @(15, 44427, "Especial", 65,  24000, 24000, 24000, "$/bandeja 7 kilos",            3429, 7),
    @(16, 44413, "Primera",  35,  3500,  3500,  3500,  "$/kilo (en caja de 15 kilos)", 3500, 1),
    @(17, 44421, "Segunda",  50,  3200,  3200,  3200,  "$/kilo (en caja de 15 kilos)", 3200, 1),
    @(18, 44377, "Segunda",  40,  3500,  3500,  3500,  "$/kilo (en caja de 15 kilos)", 3500, 1),
    @(19, 44426, "Especial", 30,  4500,  4500,  4500,  "$/kilo (en caja de 15 kilos)", 4500, 1),
    @(20, 44426, "Primera",  45,  3500,  3500,  3500,  "$/kilo (en caja de 15 kilos)", 3500, 1),
    @(21, 44354, "Tercera",  95,  3500,  3500,  3500,  "$/kilo (en caja de 15 kilos)", 3500, 1),
    @(22, 44412, "Primera",  65,  3200,  3200,  3200,  "$/kilo (en caja de 15 kilos)", 3200, 1),
    @(23, 44435, "Especial", 130, 3500,  4500,  3885,  "$/kilo (en caja de 15 kilos)", 3885, 1),
    @(24, 44435, "Primera",  80,  21000, 21000, 21000, "$/bandeja 7 kilos",            3000, 7),
    @(25, 44431, "Primera",  80,  21000, 21000, 21000, "$/bandeja 7 kilos",            3000, 7),
    @(26, 44405, "Segunda",  50,  3200,  3200,  3200,  "$/kilo (en caja de 15 kilos)", 3200, 1),
    @(27, 44417, "Segunda",  50,  3200,  3200,  3200,  "$/kilo (en caja de 15 kilos)", 3200, 1),
    @(28, 44419, "Segunda",  70,  3200,  3200,  3200,  "$/kilo (en caja de 15 kilos)", 3200, 1)
)

foreach ($row in $existingRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]     # D Fecha
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 12).Value = $row[2]    # L Calidad
    $ws.Cells.Item($r, 13).Value = $row[3]    # M Volumen
    $ws.Cells.Item($r, 14).Value = $row[4]    # N Precio minimo
    $ws.Cells.Item($r, 15).Value = $row[5]    # O Precio maximo
    $ws.Cells.Item($r, 16).Value = $row[6]    # P Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $row[7]    # Q Unidad de comercializacion
    $ws.Cells.Item($r, 19).Value = $row[8]    # S Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $row[9]    # T Kg / unidad
}

# --- Part 2: rows 29-36 are brand-new / fully-rewritten rows. All columns A-T are set.
$newRows = @(
    @(29, 44420, "Primera",  35,  3500,  3500,  3500,  "$/kilo (en caja de 15 kilos)", 3500, 1),
    @(30, 44420, "Segunda",  40,  3200,  3200,  3200,  "$/kilo (en caja de 15 kilos)", 3200, 1),
    @(31, 44161, "Primera",  65,  2300,  2300,  2300,  "$/kilo (en caja de 15 kilos)", 2300, 1),
    @(32, 44161, "Segunda",  55,  2000,  2000,  2000,  "$/kilo (en caja de 15 kilos)", 2000, 1),
    @(33, 44159, "Primera",  120, 2300,  2500,  2408,  "$/kilo (en caja de 15 kilos)", 2408, 1),
    @(34, 44433, "Especial", 20,  4500,  4500,  4500,  "$/kilo (en caja de 15 kilos)", 4500, 1),
    @(35, 44160, "Primera",  120, 2200,  2300,  2246,  "$/kilo (en caja de 15 kilos)", 2246, 1),
    @(36, 44432, "Especial", 30,  4500,  4500,  4500,  "$/kilo (en caja de 15 kilos)", 4500, 1)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = 10                            # A Mercado ID
    $ws.Cells.Item($r, 2).Value = "Vega Modelo de Temuco"       # B Mercado
    $ws.Cells.Item($r, 3).Value = "La Araucanía"                 # C Región
    $ws.Cells.Item($r, 4).Value = $row[1]                        # D Fecha
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = 9                              # E Codreg
    $ws.Cells.Item($r, 6).Value = "Fruta"                        # F Tipo
    $ws.Cells.Item($r, 7).Value = 100107                         # G Producto ID
    $ws.Cells.Item($r, 8).Value = "Otros"                        # H Producto
    $ws.Cells.Item($r, 9).Value = 100107002                      # I Categoría ID
    $ws.Cells.Item($r, 10).Value = "Chirimoya"                   # J Categoría
    $ws.Cells.Item($r, 11).Value = "Cultivar IV Región"          # K Variedad
    $ws.Cells.Item($r, 12).Value = $row[2]                       # L Calidad
    $ws.Cells.Item($r, 13).Value = $row[3]                       # M Volumen
    $ws.Cells.Item($r, 14).Value = $row[4]                       # N Precio minimo
    $ws.Cells.Item($r, 15).Value = $row[5]                       # O Precio maximo
    $ws.Cells.Item($r, 16).Value = $row[6]                       # P Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $row[7]                       # Q Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = "Provincia del Elquí"         # R Origen
    $ws.Cells.Item($r, 19).Value = $row[8]                       # S Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $row[9]                       # T Kg / unidad
}
